$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename column headers: "_old" -> "_FV2404" (cols A-J), "_new" -> "_FV2410" (cols L-U)
# ---------------------------------------------------------------------------
$headersFV2404 = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")
$headersFV2410 = @("Segmentname_FV2410","Segmentgruppe_FV2410","Segment_FV2410","Datenelement_FV2410","Segment ID_FV2410","Code_FV2410","Qualifier_FV2410","Beschreibung_FV2410","Bedingungsausdruck_FV2410","Bedingung_FV2410")

for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersFV2404[$i]
}
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersFV2410[$i]
}

# ---------------------------------------------------------------------------
# 2. Turn A1:U72 into an Excel Table ("Table1"), while avoiding an automatic
#    header-row dxf override (the engine bakes one in whenever the header
#    range already carries non-default cell formatting at Add()-time).
#    Trick: stash the header's format, build the table against a fresh
#    (unstyled) header row holding the final header text, then restore the
#    original formatting afterwards.
# ---------------------------------------------------------------------------

# 2a. Stash the current header formatting (format only) into a scratch row
#     well beyond the used range.
$ws.Range("A1:U1").Copy() | Out-Null
$ws.Range("A74:U74").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 2b. Insert a fresh, unstyled row above row 1 and copy the (already renamed)
#     header text into it.
$ws.Rows.Item(1).Insert()
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(1, $c).Value = $ws.Cells.Item(2, $c).Value2
}

# 2c. Build the table against this clean header row.
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U73"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# 2d. Drop the now-duplicate, still-styled header row (former row 1, now row 2).
$ws.Rows.Item(2).Delete()

# 2e. Restore the original header formatting from the scratch row, then
#     remove the scratch row entirely so the used range stays A1:U72.
$ws.Range("A74:U74").Copy() | Out-Null
$ws.Range("A1:U1").PasteSpecial(-4122) | Out-Null     # xlPasteFormats
$ws.Rows.Item(74).Delete()

# ---------------------------------------------------------------------------
# 3. Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
